# Journal mis à jour
# Updates the "06.01" week block (rows 27-34) of the journal sheet:
#  - row 28 (Date row): hours 1 -> 2
#  - row 29: new task text + 4 hours
#  - row 30: new task text + 1 hour
#  - row 33: total (formula SUM(D28:D32)) recalculates automatically
#  - row 34: new "réflexion personnelle" text
#  - D84 grand-total formula recalculates automatically
# Also updates the view state (scrolled-to row / current selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits -----------------------------------------------------------

$ws.Range("D28").Value = 2

$ws.Range("B29").Value = "Mise en place du contrôle du rover avec la manette"
$ws.Range("D29").Value = 4

$ws.Range("B30").Value = "Configuration de toutes les options de la vue sur la manette"
$ws.Range("D30").Value = 1

$ws.Range("B34").Value = "Aujourd'hui il y a eu beaucoup d'avancements et de réussites, on a réussi à mettre tout ce qu'on avait prévu en place, et même plus que ça. Je pense que cette journée s'est très très bien passée."

# --- View state -------------------------------------------------------------
# Scroll the frozen (bottom) pane so that row 20 is the first visible row,
# and restore the selection to the newly-filled-in comment cell B34:D34.

$win = $excel.ActiveWindow
$win.SplitRow = 19
$ws.Range("B34:D34").Select() | Out-Null
